$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 8.678923347133891
$ws.Range("C2").Value = 6.242821631731166
$ws.Range("D2").Value = 4.609631130861245
$ws.Range("F2").Value = 19.44384730135743
$ws.Range("G2").Value = 20.33466993852254
$ws.Range("H2").Value = 12.38673335610381
$ws.Range("I2").Value = 17.49699073661506
$ws.Range("K2").Value = 9.110013661397437
$ws.Range("O2").Value = 17.6977912847671
# Row 3
$ws.Range("B3").Value = 8.253335560999284
$ws.Range("C3").Value = 6.084638369732843
$ws.Range("D3").Value = 4.479888746756192
$ws.Range("F3").Value = 19.54790552534374
$ws.Range("G3").Value = 20.51844023474209
$ws.Range("H3").Value = 12.44244137500632
$ws.Range("I3").Value = 17.60918362277
$ws.Range("K3").Value = 8.801715088312658
$ws.Range("O3").Value = 17.80569547209395
# Row 4
$ws.Range("B4").Value = 7.9802634063928
$ws.Range("C4").Value = 5.985018644691126
$ws.Range("D4").Value = 4.397498209306306
$ws.Range("F4").Value = 19.61757639548159
$ws.Range("G4").Value = 20.6392528311684
$ws.Range("H4").Value = 12.47857748138293
$ws.Range("I4").Value = 17.68165098023638
$ws.Range("K4").Value = 8.605587141057796
$ws.Range("O4").Value = 17.87591744838307
# Row 5
$ws.Range("B5").Value = 7.866142909426904
$ws.Range("C5").Value = 5.943840244741648
$ws.Range("D5").Value = 4.363268052732776
$ws.Range("F5").Value = 19.6474166387054
$ws.Range("G5").Value = 20.69048370412338
$ws.Range("H5").Value = 12.49378970485923
$ws.Range("I5").Value = 17.71208452441557
$ws.Range("K5").Value = 8.524029457145437
$ws.Range("O5").Value = 17.90553191869763
# Row 6
$ws.Range("B6").Value = 7.847025445206893
$ws.Range("C6").Value = 5.93696878175713
$ws.Range("D6").Value = 4.357545527165525
$ws.Range("F6").Value = 19.65245896561333
$ws.Range("G6").Value = 20.69911106871195
$ws.Range("H6").Value = 12.49634508990375
$ws.Range("I6").Value = 17.7171925644794
$ws.Range("K6").Value = 8.510390802109026
$ws.Range("O6").Value = 17.9105097001178
# Row 7
$ws.Range("B7").Value = 7.978735674821181
$ws.Range("C7").Value = 5.98446559590481
$ws.Range("D7").Value = 4.397039180491772
$ws.Range("F7").Value = 19.6179729710134
$ws.Range("G7").Value = 20.63993566493239
$ws.Range("H7").Value = 12.47878066790836
$ws.Range("I7").Value = 17.6820577605413
$ws.Range("K7").Value = 8.604493724195768
$ws.Range("O7").Value = 17.87631279602788
# Row 8
$ws.Range("B8").Value = 8.534684865487211
$ws.Range("C8").Value = 6.188823717231815
$ws.Range("D8").Value = 4.565479805233194
$ws.Range("F8").Value = 19.47852463685245
$ws.Range("G8").Value = 20.3963740077643
$ws.Range("H8").Value = 12.40554118799672
$ws.Range("I8").Value = 17.53493302561625
$ws.Range("K8").Value = 9.005177782536233
$ws.Range("O8").Value = 17.73417334865431
# Row 9
$ws.Range("B9").Value = 9.527430791145902
$ws.Range("C9").Value = 6.56796155626095
$ws.Range("D9").Value = 4.872926716794391
$ws.Range("F9").Value = 19.25110829435704
$ws.Range("G9").Value = 19.98236895495086
$ws.Range("H9").Value = 12.27720216076382
$ws.Range("I9").Value = 17.27473139334875
$ws.Range("K9").Value = 9.73346501172376
$ws.Range("O9").Value = 17.48689970813047
# Row 10
$ws.Range("B10").Value = 10.19289934999201
$ws.Range("C10").Value = 6.831124475146631
$ws.Range("D10").Value = 5.083414919929545
$ws.Range("F10").Value = 19.11234978400127
$ws.Range("G10").Value = 19.71744202342785
$ws.Range("H10").Value = 12.1921731105756
$ws.Range("I10").Value = 17.10067836147254
$ws.Range("K10").Value = 10.22973138872589
$ws.Range("O10").Value = 17.32437147368523
# Row 11
$ws.Range("B11").Value = 10.48107541979091
$ws.Range("C11").Value = 6.947093521408311
$ws.Range("D11").Value = 5.175568946027497
$ws.Range("F11").Value = 19.05542900577339
$ws.Range("G11").Value = 19.60554532532619
$ws.Range("H11").Value = 12.15549067171331
$ws.Range("I11").Value = 17.0251840898772
$ws.Range("K11").Value = 10.44641870060795
$ws.Range("O11").Value = 17.25458349192579
# Row 12
$ws.Range("B12").Value = 10.58806426946569
$ws.Range("C12").Value = 6.990439219169758
$ws.Range("D12").Value = 5.209929055513552
$ws.Range("F12").Value = 19.03477068558838
$ws.Range("G12").Value = 19.56442174964505
$ws.Range("H12").Value = 12.14188643214181
$ws.Range("I12").Value = 16.99712385090075
$ws.Range("K12").Value = 10.52712320400891
$ws.Range("O12").Value = 17.22875266681296
# Row 13
$ws.Range("B13").Value = 10.56511798514412
$ws.Range("C13").Value = 6.981129725089567
$ws.Range("D13").Value = 5.202553141717395
$ws.Range("F13").Value = 19.03917988674763
$ws.Range("G13").Value = 19.57322272720277
$ws.Range("H13").Value = 12.14480361212378
$ws.Range("I13").Value = 17.0031436804518
$ws.Range("K13").Value = 10.50980277653305
$ws.Range("O13").Value = 17.23428927206997
# Row 14
$ws.Range("B14").Value = 10.48992051680834
$ws.Range("C14").Value = 6.950671172923789
$ws.Range("D14").Value = 5.178406630369749
$ws.Range("F14").Value = 19.0537114431356
$ws.Range("G14").Value = 19.60213696976492
$ws.Range("H14").Value = 12.15436570249964
$ws.Range("I14").Value = 17.02286498968176
$ws.Range("K14").Value = 10.45308561546245
$ws.Range("O14").Value = 17.25244642177617
# Row 15
$ws.Range("B15").Value = 10.44358038204577
$ws.Range("C15").Value = 6.931939431594762
$ws.Range("D15").Value = 5.16354577590799
$ws.Range("F15").Value = 19.06272930698464
$ws.Range("G15").Value = 19.62001076516056
$ws.Range("H15").Value = 12.16026006387713
$ws.Range("I15").Value = 17.03501354205869
$ws.Range("K15").Value = 10.41816759240099
$ws.Range("O15").Value = 17.26364586758686
# Row 16
$ws.Range("B16").Value = 10.17376955365134
$ws.Range("C16").Value = 6.823467573490409
$ws.Range("D16").Value = 5.077318344470032
$ws.Range("F16").Value = 19.11619487130429
$ws.Range("G16").Value = 19.72492901622777
$ws.Range("H16").Value = 12.19461053598373
$ws.Range("I16").Value = 17.10568604838816
$ws.Range("K16").Value = 10.21538380524018
$ws.Range("O16").Value = 17.32901573161573
# Row 17
$ws.Range("B17").Value = 10.00448814542142
$ws.Range("C17").Value = 6.755941632624676
$ws.Range("D17").Value = 5.023484711311172
$ws.Range("F17").Value = 19.15058584829779
$ws.Range("G17").Value = 19.79150756410686
$ws.Range("H17").Value = 12.21619461972565
$ws.Range("I17").Value = 17.14998340103729
$ws.Range("K17").Value = 10.08862529649446
$ws.Range("O17").Value = 17.37018016520381
# Row 18
$ws.Range("B18").Value = 9.905754630423539
$ws.Range("C18").Value = 6.716751645247658
$ws.Range("D18").Value = 4.992183370299171
$ws.Range("F18").Value = 19.17095004055009
$ws.Range("G18").Value = 19.83061252722142
$ws.Range("H18").Value = 12.22879728552724
$ws.Range("I18").Value = 17.17580886771231
$ws.Range("K18").Value = 10.0148661251618
$ws.Range("O18").Value = 17.39424718022493
# Row 19
$ws.Range("B19").Value = 9.872091836258729
$ws.Range("C19").Value = 6.7034233123594
$ws.Range("D19").Value = 4.981527900339097
$ws.Range("F19").Value = 19.17794506533815
$ws.Range("G19").Value = 19.8439917497744
$ws.Range("H19").Value = 12.23309665363767
$ws.Range("I19").Value = 17.18461254739705
$ws.Range("K19").Value = 9.989747887758387
$ws.Range("O19").Value = 17.40246289577338
# Row 20
$ws.Range("B20").Value = 10.02265030725486
$ws.Range("C20").Value = 6.763166427709444
$ws.Range("D20").Value = 5.029250485637374
$ws.Range("F20").Value = 19.14686445537282
$ws.Range("G20").Value = 19.78433618551354
$ws.Range("H20").Value = 12.21387749731238
$ws.Range("I20").Value = 17.14523199314816
$ws.Range("K20").Value = 10.10220739030973
$ws.Range("O20").Value = 17.36575774403538
# Row 21
$ws.Range("B21").Value = 10.51206615521509
$ws.Range("C21").Value = 6.959633268576608
$ws.Range("D21").Value = 5.1855137527771
$ws.Range("F21").Value = 19.04941880989638
$ws.Range("G21").Value = 19.5936101716888
$ws.Range("H21").Value = 12.15154931161361
$ws.Range("I21").Value = 17.0170580533134
$ws.Range("K21").Value = 10.46978180775909
$ws.Range("O21").Value = 17.24709704079514
# Row 22
$ws.Range("B22").Value = 10.81945546440679
$ws.Range("C22").Value = 7.084705310971352
$ws.Range("D22").Value = 5.284505501196467
$ws.Range("F22").Value = 18.99095909823779
$ws.Range("G22").Value = 19.47624616622089
$ws.Range("H22").Value = 12.11248449160791
$ws.Range("I22").Value = 16.9363645760312
$ws.Range("K22").Value = 10.70212499241898
$ws.Range("O22").Value = 17.1730215751014
# Row 23
$ws.Range("B23").Value = 10.65654995774345
$ws.Range("C23").Value = 7.018266169020938
$ws.Range("D23").Value = 5.231964450265273
$ws.Range("F23").Value = 19.021680404833
$ws.Range("G23").Value = 19.53821553381036
$ws.Range("H23").Value = 12.13318150733976
$ws.Range("I23").Value = 16.97915139069462
$ws.Range("K23").Value = 10.57885465870405
$ws.Range("O23").Value = 17.21223894364193
# Row 24
$ws.Range("B24").Value = 10.01444358651902
$ws.Range("C24").Value = 6.759901243061968
$ws.Range("D24").Value = 5.02664487295829
$ws.Range("F24").Value = 19.14854505232894
$ws.Range("G24").Value = 19.78757578716164
$ws.Range("H24").Value = 12.21492446498256
$ws.Range("I24").Value = 17.14737898835401
$ws.Range("K24").Value = 10.09606967563678
$ws.Range("O24").Value = 17.36775587146542
# Row 25
$ws.Range("B25").Value = 9.269817749759989
$ws.Range("C25").Value = 6.46794861913327
$ws.Range("D25").Value = 4.792361498962611
$ws.Range("F25").Value = 19.30767470264172
$ws.Range("G25").Value = 20.08751116158841
$ws.Range("H25").Value = 12.31029088358539
$ws.Range("I25").Value = 17.34210679401476
$ws.Range("K25").Value = 9.543014845901734
$ws.Range("O25").Value = 17.87631279602788
